$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Combine the separate "Motor PWM" / "Servo PWM" threads into a single
# "PWM" row. Row 4 currently holds "Servo PWM" -> rename it to "PWM".
$ws.Range("A4").Value = "PWM"

# Row 8 currently holds "Motor PWM" -> remove it entirely; Excel shifts
# every row below it up by one, closing the gap.
$ws.Rows.Item(8).Delete()

# Restore the selection to match the post-edit layout.
$ws.Range("A9:E13").Select()
